# Insert a new weekly price record as row 134 on the "Albahaca" sheet.
# This pushes the existing rows 134-166 down to 135-167 (dimension grows
# from A1:R166 to A1:R167) and populates the new row with the new record's
# data, matching column order:
#   A Mercado ID | B Mercado | C Región | D Fecha | E Codreg | F Categoría ID
#   G Categoría | H Variedad | I Calidad | J Volumen | K Precio mínimo
#   L Precio máximo | M Precio promedio ponderado | N Unidad de comercialización
#   O Origen | P Precio $/Kg | Q Kg o Unidades | R Clasificación

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134..166 down to 135..167, leaving a blank row 134 (it
# inherits row 133's formatting, e.g. the date style on column D).
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new record.
$ws.Cells.Item(134, 1).Value  = 8
$ws.Cells.Item(134, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(134, 3).Value  = "Coquimbo"
$ws.Cells.Item(134, 4).Value  = 44995
$ws.Cells.Item(134, 5).Value  = 4
$ws.Cells.Item(134, 6).Value  = 100112052
$ws.Cells.Item(134, 7).Value  = "Albahaca"
$ws.Cells.Item(134, 8).Value  = "Sin especificar"
$ws.Cells.Item(134, 9).Value  = "Primera"
$ws.Cells.Item(134, 10).Value = 1000
$ws.Cells.Item(134, 11).Value = 2800
$ws.Cells.Item(134, 12).Value = 3000
$ws.Cells.Item(134, 13).Value = 2900
$ws.Cells.Item(134, 14).Value = "`$/docena de matas"
$ws.Cells.Item(134, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(134, 16).Value = 483
$ws.Cells.Item(134, 17).Value = 6
$ws.Cells.Item(134, 18).Value = "Hortaliza"
